$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Contact")

# New column L: header (row1) mirrors the "...Message" shared string,
# data (row2) mirrors the "No X Found" shared string - same pattern as
# the existing "no Organization Found Message" / "No Organization Found"
# pair in column K.
$ws.Range("L2").Value = "No Contact Found !"
$ws.Range("L1").Value = "no Contact Found Message"

# Give L1 the same header style (fill/border) as the other header cells.
$ws.Range("K1").Copy()
$ws.Range("L1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Match the column width recorded in the workbook (stored width "25").
$ws.Columns.Item(12).ColumnWidth = 24.166666666666668

# Reproduce the new selection: whole column L selected, which is what
# produces activeCell="L1" sqref="L1:L1048576" in the saved view state.
$ws.Activate()
$ws.Columns.Item(12).Select()
